$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.559.07'
$ws.Range("E2").Value = '  +3.52%  '

$ws.Range("D3").Value = '2.403.09'
$ws.Range("E3").Value = '  +3.47%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '553.43'
$ws.Range("E5").Value = '  +3.28%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '136.41'
$ws.Range("E6").Value = '  +2.05%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.11%  '

$ws.Range("E8").Value = '  +2.76%  '

$ws.Range("E9").Value = '  +7.34%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '5.87'
$ws.Range("E10").Value = '  +7.84%  '

$ws.Range("B11").Value = 'Cardano'
$ws.Range("C11").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.363'
$ws.Range("E11").Value = '  +2.11%  '

$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.151'
$ws.Range("E12").Value = '  -0.88%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '24.67'
$ws.Range("E13").Value = '  +5.00%  '

$ws.Range("D14").Value = '2.829.91'
$ws.Range("E14").Value = '  +3.57%  '

$ws.Range("D15").Value = '59.446.24'
$ws.Range("E15").Value = '  +3.37%  '

$ws.Range("E16").Value = '  +5.84%  '

$ws.Range("D17").Value = '2.382.01'
$ws.Range("E17").Value = '  +2.35%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.28'
$ws.Range("E18").Value = '  +6.95%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.41'
$ws.Range("E19").Value = '  +4.74%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '336.56'
$ws.Range("E20").Value = '  +1.85%  '

$ws.Range("E21").Value = '  +5.87%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.999'
$ws.Range("E22").Value = '  -0.03%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '64.64'
$ws.Range("E23").Value = '  +4.24%  '

$ws.Range("E24").Value = '  +1.36%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.00'
$ws.Range("E25").Value = '  -0.08%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.45'
$ws.Range("E26").Value = '  +0.51%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.38'
$ws.Range("E27").Value = '  -1.16%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.80'
$ws.Range("E28").Value = '  +3.00%  '

$ws.Range("D29").Value = '0.0₃0765'
$ws.Range("E29").Value = '  +5.61%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '171.36'
$ws.Range("E30").Value = '  +0.82%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.27'
$ws.Range("E31").Value = '  +3.49%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '18.74'
$ws.Range("E32").Value = '  +2.02%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.03'
$ws.Range("E33").Value = '  +0.01%  '

$ws.Range("E34").Value = '  -0.03%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.32'
$ws.Range("E35").Value = '  +4.12%  '

$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.29'
$ws.Range("E36").Value = '  +4.87%  '

$ws.Range("B37").Value = 'FirstDigitalUSD'
$ws.Range("C37").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.999'
$ws.Range("E37").Value = '  -0.11%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.66'
$ws.Range("E38").Value = '  +3.06%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '40.18'
$ws.Range("E39").Value = '  +3.13%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.423'
$ws.Range("E40").Value = '  +13.26%  '

$ws.Range("E41").Value = '  +4.26%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '296.22'
$ws.Range("E42").Value = '  +5.28%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '142.18'
$ws.Range("E43").Value = '  -1.20%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0964'
$ws.Range("E44").Value = '  +3.19%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0526'
$ws.Range("E45").Value = '  +5.45%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '19.21'
$ws.Range("E46").Value = '  +1.89%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.571'
$ws.Range("E47").Value = '  +2.62%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0227'
$ws.Range("E48").Value = '  +5.79%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.400'
$ws.Range("E49").Value = '  +4.18%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '11.04'
$ws.Range("E50").Value = '  -0.20%  '

$ws.Range("E51").Value = '  +5.44%  '
